# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'76.874.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.77%  '

# Row 3
$ws.Range("D3").Value = "'3.134.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.26%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = "'200.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.33%  '

# Row 6
$ws.Range("D6").Value = "'626.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").Value = "'0.225"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +13.93%  '

# Row 9
$ws.Range("D9").Value = "'0.571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.12%  '

# Row 10
$ws.Range("D10").Value = "'3.131.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.08%  '

# Row 11
$ws.Range("D11").Value = "'0.519"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +19.56%  '

# Row 12
$ws.Range("E12").Value = '  +0.95%  '

# Row 13
$ws.Range("D13").Value = "'5.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.76%  '

# Row 14
$ws.Range("D14").Value = "'3.705.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.95%  '

# Row 15
$ws.Range("D15").Value = "'0.0000222"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +19.35%  '

# Row 16
$ws.Range("D16").Value = "'30.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.83%  '

# Row 17
$ws.Range("D17").Value = "'76.676.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.64%  '

# Row 18
$ws.Range("D18").Value = "'3.107.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.87%  '

# Row 19
$ws.Range("D19").Value = "'13.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.64%  '

# Row 20
$ws.Range("D20").Value = "'9.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.99%  '

# Row 21
$ws.Range("D21").Value = "'2.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +24.46%  '

# Row 22
$ws.Range("D22").Value = "'400.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.30%  '

# Row 23
$ws.Range("D23").Value = "'4.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.07%  '

# Row 24
$ws.Range("D24").Value = "'6.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.85%  '

# Row 25
$ws.Range("D25").Value = "'3.289.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.01%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = "'4.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.59%  '

# Row 27
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = "'74.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.80%  '

# Row 28
$ws.Range("D28").Value = "'10.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.87%  '

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.0000114"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.99%  '

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.27%  '

# Row 32
$ws.Range("D32").Value = "'8.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.33%  '

# Row 33
$ws.Range("D33").Value = "'1.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.05%  '

# Row 34
$ws.Range("D34").Value = "'517.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.01%  '

# Row 35
$ws.Range("D35").Value = "'1.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.37%  '

# Row 36
$ws.Range("D36").Value = "'0.134"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +20.87%  '

# Row 37
$ws.Range("D37").Value = "'21.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.65%  '

# Row 38
$ws.Range("E38").Value = '  -0.14%  '

# Row 39
$ws.Range("D39").Value = "'163.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.09%  '

# Row 40
$ws.Range("D40").Value = "'0.387"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.65%  '

# Row 41
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = "'20.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.32%  '

# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'193.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.46%  '

# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.04%  '

# Row 44
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = "'0.103"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.66%  '

# Row 45
$ws.Range("D45").Value = "'5.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.85%  '

# Row 46
$ws.Range("D46").Value = "'0.793"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +18.42%  '

# Row 47
$ws.Range("D47").Value = "'1.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.57%  '

# Row 48
$ws.Range("E48").Value = '  +6.32%  '

# Row 49
$ws.Range("D49").Value = "'42.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.76%  '

# Row 50
$ws.Range("D50").Value = "'2.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.94%  '

# Row 51
$ws.Range("D51").Value = "'0.615"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.59%  '
